# Language workbook update: add a new "criteria" / "CRITERIA" key/value
# pair as a new row 11 (pushing the season/atmosphere/climate/region rows
# down by one), and leave the selection sitting on the newly added row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("en")

# Insert a fresh row above the current row 11 ("season_winter"/"Winter"),
# shifting everything from row 11 down to row 12 onward.
$ws.Range("A11").EntireRow.Insert()

# Populate the new row with the criteria key/value pair.
$ws.Range("A11").Value = "criteria"
$ws.Range("B11").Value = "CRITERIA"

# Match the saved selection state (active cell on the newly inserted row).
[void]$ws.Range("A11").Select()
